# DS - Address Change - Login, Guest, Email and Mobile Changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the address used by the Linking tests (UNIT 35 -> UNIT 66) ---
$ws.Range("M4").Value = "UNIT 66, 653-659 GEORGE ST, HAYMARKET NSW 2000"
$ws.Range("N4").Value = "UNIT 66, 653-659 GEORGE ST, HAYMARKET NSW 2000"
$ws.Range("M5").Value = "UNIT 66, 653-659 GEORGE ST, HAYMARKET NSW 2000"
$ws.Range("N5").Value = "UNIT 66, 653-659 GEORGE ST, HAYMARKET NSW 2000"

# --- CheckDemeritsTest now uses the sa021 mailbox ---
$ws.Range("B7").Value = "sa021@mailinator.com"

# --- New header columns: Email Id / Mobile Number ---
$ws.Range("O1").Value = "Email Id"
$ws.Range("O1").Style = $ws.Range("A1").Style

$ws.Range("P1").Value = "Mobile Number"
$ws.Range("P1").Interior.ColorIndex = 6
$ws.Range("P1").NumberFormat = "@"

# --- New row 8: ChangeAddressDetailsTest ---
$ws.Range("A8").Value = "ChangeAddressDetailsTest"
$ws.Range("B8").Value = "sa021@mailinator.com"
$ws.Range("B8").Style = $ws.Range("B7").Style
$ws.Range("C8").Value = "Pa$$w0rd"

# --- New row 9: ChangeAddressDetailsGuestTest ---
$ws.Range("A9").Value = "ChangeAddressDetailsGuestTest"
$ws.Range("B9").Style = $ws.Range("B7").Style
$ws.Range("F9").Value = 45164185
$ws.Range("G9").Value = 98951626
$ws.Range("H9").Value = "SNSWSN"

# --- New row 10: ChangeMobileAndEmailTest ---
$ws.Range("A10").Value = "ChangeMobileAndEmailTest"
$ws.Range("B10").Value = "sa021@mailinator.com"
$ws.Range("B10").Style = $ws.Range("B7").Style
$ws.Range("C10").Value = "Pa$$w0rd"
$ws.Range("O10").Value = "sa021@mailinator4.com"
$ws.Range("O10").Style = $ws.Range("B7").Style
$ws.Range("P10").Value = "0451764467"
$ws.Range("P10").NumberFormat = "@"

# --- Column sizing for the new / widened columns ---
$ws.Columns.Item(1).ColumnWidth = 30.42578125
$ws.Range("M1:N1").EntireColumn.ColumnWidth = 52.7109375
$ws.Columns.Item(15).ColumnWidth = 21.5703125
$ws.Columns.Item(16).ColumnWidth = 15.140625

# --- Selection mirrors the saved workbook state ---
$ws.Range("D14").Select()
